$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Fill in the previously-empty row 4 values
$ws.Range("A4").Value = 44580
$ws.Range("B4").Value = 101
$ws.Range("C4").Value = 224
$ws.Range("D4").Value = 177
$ws.Range("E4").Value = 4368
$ws.Range("F4").Value = 1525
$ws.Range("G4").Value = 3197
$ws.Range("H4").Value = 57422
$ws.Range("I4").Value = 40634
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 203
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 135
$ws.Range("O4").Value = 125
$ws.Range("P4").Value = 122
$ws.Range("Q4").Value = 2
$ws.Range("R4").Value = 543

# Add a new row 5 with the next data point
$ws.Range("A5").Value = 44588
$ws.Range("B5").Value = 102
$ws.Range("C5").Value = 224
$ws.Range("D5").Value = 177
$ws.Range("E5").Value = 4368
$ws.Range("F5").Value = 1531
$ws.Range("G5").Value = 3207
$ws.Range("H5").Value = 57794
$ws.Range("I5").Value = 40652
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 203
$ws.Range("M5").Value = 0
$ws.Range("N5").Value = 135
$ws.Range("O5").Value = 125
$ws.Range("P5").Value = 122
$ws.Range("Q5").Value = 3
$ws.Range("R5").Value = 550

# Copy styles from row 4 (Date col date-format, rest numeric) down to row 5
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B4:R4").Copy()
$ws.Range("B5:R5").PasteSpecial(-4122)

# Adjust the view to match target (scroll to column F, select F6)
$ws.Application.ActiveWindow.ScrollColumn = 6
$ws.Range("F6").Select()
